$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1558.4286
$ws.Range("J17").Value = 1558.4286
$ws.Range("L17").Value = 4675.2858
$ws.Range("N17").Value = -5011.2858
$ws.Range("H28").Value = 2082.8125
$ws.Range("I28").Value = 1688.6333
$ws.Range("K28").Value = 1688.6333
$ws.Range("M28").Value = -1203.6333
$ws.Range("H33").Value = 189.23077
$ws.Range("I33").Value = 191.75
$ws.Range("J33").Value = 159
$ws.Range("K33").Value = 191.75
$ws.Range("L33").Value = 159
$ws.Range("M33").Value = 37.25
$ws.Range("N33").Value = -617
$ws.Range("H38").Value = 583.5714
$ws.Range("I38").Value = 149.33333
$ws.Range("J38").Value = 909.25
$ws.Range("K38").Value = 447.99999
$ws.Range("L38").Value = 2727.75
$ws.Range("M38").Value = -75.99998999999997
$ws.Range("N38").Value = -3471.75
$ws.Range("H64").Value = 8110.5
$ws.Range("J64").Value = 7998.5
$ws.Range("L64").Value = 7998.5
$ws.Range("N64").Value = -8494.5
$ws.Range("H67").Value = 8110.5
$ws.Range("J67").Value = 7998.5
$ws.Range("L67").Value = 7998.5
$ws.Range("N67").Value = -9714.5
$ws.Range("H70").Value = 9780.538
$ws.Range("I70").Value = 3200
$ws.Range("J70").Value = 10977
$ws.Range("K70").Value = 9600
$ws.Range("L70").Value = 32931
$ws.Range("M70").Value = -9330
$ws.Range("N70").Value = -33471
$ws.Range("H73").Value = 9780.538
$ws.Range("I73").Value = 3200
$ws.Range("J73").Value = 10977
$ws.Range("K73").Value = 9600
$ws.Range("L73").Value = 32931
$ws.Range("M73").Value = -8664
$ws.Range("N73").Value = -34803
$ws.Range("H88").Value = 1504.5333
$ws.Range("I88").Value = 731.4
$ws.Range("J88").Value = 1891.1
$ws.Range("K88").Value = 731.4
$ws.Range("L88").Value = 1891.1
$ws.Range("M88").Value = -325.4
$ws.Range("N88").Value = -2703.1
$ws.Range("H91").Value = 1504.5333
$ws.Range("I91").Value = 731.4
$ws.Range("J91").Value = 1891.1
$ws.Range("K91").Value = 731.4
$ws.Range("L91").Value = 1891.1
$ws.Range("M91").Value = 672.6
$ws.Range("N91").Value = -4699.1
$ws.Range("H116").Value = 7119
$ws.Range("I116").Value = 7505.5557
$ws.Range("K116").Value = 7505.5557
$ws.Range("M116").Value = -4063.5557
$ws.Range("H135").Value = 1990.5172
$ws.Range("I135").Value = 1591.625
$ws.Range("K135").Value = 14324.625
$ws.Range("M135").Value = -11789.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2841
$ws.Range("I74").Value = 2786.6538
$ws.Range("K74").Value = 2786.6538
$ws.Range("M74").Value = -1912.6538
$ws.Range("H77").Value = 2841
$ws.Range("I77").Value = 2786.6538
$ws.Range("K77").Value = 13933.269
$ws.Range("M77").Value = -9565.269
$ws.Range("H88").Value = 13296.182
$ws.Range("J88").Value = 15818.333
$ws.Range("L88").Value = 15818.333
$ws.Range("N88").Value = -16630.333
$ws.Range("H91").Value = 13296.182
$ws.Range("J91").Value = 15818.333
$ws.Range("L91").Value = 15818.333
$ws.Range("N91").Value = -18626.333
$ws.Range("H110").Value = 435.875
$ws.Range("I110").Value = 409.5
$ws.Range("J110").Value = 515
$ws.Range("K110").Value = 409.5
$ws.Range("L110").Value = 515
$ws.Range("M110").Value = 1635.5
$ws.Range("N110").Value = -4605
$ws.Range("H132").Value = 3559.625
$ws.Range("I132").Value = 3407.3076
$ws.Range("K132").Value = 10221.9228
$ws.Range("M132").Value = -7691.9228

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 6705.75
$ws.Range("I86").Value = 3716.5417
$ws.Range("J86").Value = 12684.167
$ws.Range("K86").Value = 3716.5417
$ws.Range("L86").Value = 12684.167
$ws.Range("M86").Value = -2593.5417
$ws.Range("N86").Value = -14930.167
$ws.Range("H89").Value = 6705.75
$ws.Range("I89").Value = 3716.5417
$ws.Range("J89").Value = 12684.167
$ws.Range("K89").Value = 18582.7085
$ws.Range("L89").Value = 63420.835
$ws.Range("M89").Value = -12966.7085
$ws.Range("N89").Value = -74652.83499999999
$ws.Range("H107").Value = 1409.6471
$ws.Range("I107").Value = 1190.0834
$ws.Range("K107").Value = 1190.0834
$ws.Range("M107").Value = 729.9166
$ws.Range("H134").Value = 9934.182000000001
$ws.Range("I134").Value = 9279.833000000001
$ws.Range("K134").Value = 27839.499
$ws.Range("M134").Value = -25304.499

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 10719.2
$ws.Range("I7").Value = 26447.5
$ws.Range("J7").Value = 233.66667
$ws.Range("K7").Value = 26447.5
$ws.Range("L7").Value = 233.66667
$ws.Range("M7").Value = -26334.5
$ws.Range("N7").Value = -459.66667
$ws.Range("H22").Value = 567.75
$ws.Range("I22").Value = 374.5
$ws.Range("J22").Value = 761
$ws.Range("K22").Value = 374.5
$ws.Range("L22").Value = 761
$ws.Range("M22").Value = -24.5
$ws.Range("N22").Value = -1461
$ws.Range("H31").Value = 2233.5
$ws.Range("I31").Value = 2105.3635
$ws.Range("J31").Value = 2307.6843
$ws.Range("K31").Value = 2105.3635
$ws.Range("L31").Value = 2307.6843
$ws.Range("M31").Value = -1810.3635
$ws.Range("N31").Value = -2897.6843
$ws.Range("H34").Value = 2233.5
$ws.Range("I34").Value = 2105.3635
$ws.Range("J34").Value = 2307.6843
$ws.Range("K34").Value = 2105.3635
$ws.Range("L34").Value = 2307.6843
$ws.Range("M34").Value = -1903.3635
$ws.Range("N34").Value = -2711.6843
$ws.Range("H47").Value = 20349.5
$ws.Range("J47").Value = 20349.5
$ws.Range("L47").Value = 20349.5
$ws.Range("N47").Value = -21481.5
$ws.Range("H62").Value = 4532.75
$ws.Range("I62").Value = 3166
$ws.Range("J62").Value = 5899.5
$ws.Range("K62").Value = 3166
$ws.Range("L62").Value = 5899.5
$ws.Range("M62").Value = -2542
$ws.Range("N62").Value = -7147.5
$ws.Range("H65").Value = 4532.75
$ws.Range("I65").Value = 3166
$ws.Range("J65").Value = 5899.5
$ws.Range("K65").Value = 15830
$ws.Range("L65").Value = 29497.5
$ws.Range("M65").Value = -12710
$ws.Range("N65").Value = -35737.5
$ws.Range("H96").Value = 29000
$ws.Range("J96").Value = 29000
$ws.Range("L96").Value = 29000
$ws.Range("N96").Value = -34492

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 135494.95
$ws.Range("J37").Value = 135494.95
$ws.Range("L37").Value = 406484.85
$ws.Range("N37").Value = -406708.85
$ws.Range("H99").Value = 8822
$ws.Range("I99").Value = 350.8
$ws.Range("K99").Value = 1052.4
$ws.Range("M99").Value = 1193.6
$ws.Range("H108").Value = 11504.25
$ws.Range("I108").Value = 6013.5
$ws.Range("K108").Value = 18040.5
$ws.Range("M108").Value = -15160.5
$ws.Range("H119").Value = 10999.25
$ws.Range("I119").Value = 2000
$ws.Range("K119").Value = 6000
$ws.Range("M119").Value = -1162
$ws.Range("H140").Value = 9235.043
$ws.Range("I140").Value = 1593.9546
$ws.Range("K140").Value = 4781.8638
$ws.Range("M140").Value = 398.1361999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6188.4165
$ws.Range("I70").Value = 5602.5
$ws.Range("J70").Value = 6774.3335
$ws.Range("K70").Value = 5602.5
$ws.Range("L70").Value = 6774.3335
$ws.Range("M70").Value = -5332.5
$ws.Range("N70").Value = -7314.3335
$ws.Range("H73").Value = 6188.4165
$ws.Range("I73").Value = 5602.5
$ws.Range("J73").Value = 6774.3335
$ws.Range("K73").Value = 5602.5
$ws.Range("L73").Value = 6774.3335
$ws.Range("M73").Value = -4666.5
$ws.Range("N73").Value = -8646.333500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H103").Value = 36497.668
$ws.Range("J103").Value = 36497.668
$ws.Range("L103").Value = 36497.668
$ws.Range("N103").Value = -38841.668
$ws.Range("H122").Value = 3466.25
$ws.Range("I122").Value = 3356.4285
$ws.Range("K122").Value = 10069.2855
$ws.Range("M122").Value = -7619.2855
$ws.Range("H132").Value = 2928.5762
$ws.Range("I132").Value = 2311.394
$ws.Range("J132").Value = 3711.923
$ws.Range("K132").Value = 6934.181999999999
$ws.Range("L132").Value = 11135.769
$ws.Range("M132").Value = -4404.181999999999
$ws.Range("N132").Value = -16195.769
$ws.Range("H136").Value = 2747.4866
$ws.Range("I136").Value = 2461
$ws.Range("K136").Value = 7383
$ws.Range("M136").Value = -4833

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 11470.477
$ws.Range("J14").Value = 10744.5
$ws.Range("L14").Value = 10744.5
$ws.Range("N14").Value = -11080.5
$ws.Range("H62").Value = 6399
$ws.Range("I62").Value = 5099
$ws.Range("J62").Value = 8999
$ws.Range("K62").Value = 5099
$ws.Range("L62").Value = 8999
$ws.Range("M62").Value = -4475
$ws.Range("N62").Value = -10247
$ws.Range("H65").Value = 6399
$ws.Range("I65").Value = 5099
$ws.Range("J65").Value = 8999
$ws.Range("K65").Value = 25495
$ws.Range("L65").Value = 44995
$ws.Range("M65").Value = -22375
$ws.Range("N65").Value = -51235
$ws.Range("H107").Value = 680.2692
$ws.Range("I107").Value = 424
$ws.Range("K107").Value = 1272
$ws.Range("M107").Value = 648
$ws.Range("H132").Value = 1411.5714
$ws.Range("I132").Value = 1313.5555
$ws.Range("K132").Value = 3940.6665
$ws.Range("M132").Value = -1410.6665
